$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine extent of the data
$lastRow = $ws.UsedRange.Rows.Count

# --- Header row (row 1): columns D..G swap meaning ---
# Old header: D=codeforiati:group-code, E=codeforiati:category-name,
#             F=codeforiati:group-name,  G=codeforiati:category-code
# New header: D=codeforiati:category-name, E=codeforiati:category-code,
#             F=codeforiati:group-name,    G=codeforiati:group-code
$ws.Range("D1").Value2 = "codeforiati:category-name"
$ws.Range("E1").Value2 = "codeforiati:category-code"
$ws.Range("F1").Value2 = "codeforiati:group-name"
$ws.Range("G1").Value2 = "codeforiati:group-code"

# --- Data rows: for every row, re-map the D/E/F/G columns ---
# new D (category-name) = old E
# new E (category-code) = old G
# new F (group-name)    = old F   (unchanged)
# new G (group-code)    = old D
for ($r = 2; $r -le $lastRow; $r++) {
    $oldD = $ws.Cells.Item($r, 4).Value2
    $oldE = $ws.Cells.Item($r, 5).Value2
    $oldF = $ws.Cells.Item($r, 6).Value2
    $oldG = $ws.Cells.Item($r, 7).Value2

    $ws.Cells.Item($r, 4).Value2 = $oldE
    $ws.Cells.Item($r, 5).Value2 = $oldG
    $ws.Cells.Item($r, 6).Value2 = $oldF
    $ws.Cells.Item($r, 7).Value2 = $oldD
}
